# "added single info line in Logger" - the logger now also records the
# list of features that were added (added_features) at every call, and
# the header that used to read n_features is renamed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header rename: n_features -> added_features ---
$ws.Range("C1").Value = "added_features"

# --- fill in the new "added_features" column for the rows that already existed ---
$ws.Range("E2").Value = "['age', 'sex']"
$ws.Range("E3").Value = "['age', 'sex', 'ponytail']"

# --- new log rows appended by the extra logger call ---
$rows = @(
    @(0, $false, 10, $false, "['age', 'sex']"),
    @(0, $false, 10, 12,     "['age', 'sex']"),
    @(0, $false, 0,  $false, "['age', 'sex']"),
    @(0, $false, 0,  $false, "['age', 'sex', 'ponytail']"),
    @(0, $false, 10, $false, "['age', 'sex']"),
    @(0, $false, 10, 12,     "['age', 'sex']")
)

$r = 4
foreach ($row in $rows) {
    # column A keeps the same bordered/centered style used by the rows above it
    $ws.Range("A2").Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = 0.66
    $ws.Cells.Item($r, 7).Value = 100

    $r++
}
